# "clean cells copy and paste pev"
# Duplicates the existing "created_sheet_20230830_140614" worksheet twice,
# naming the copies after their creation timestamps, and on the second
# copy pastes column A (A2:A5) into column C (C2:C5).

$wb = $excel.ActiveWorkbook
$source = $wb.Worksheets.Item("created_sheet_20230830_140614")

# --- First new sheet: created_sheet_20230901_143340 -----------------------
# Plain duplicate of the source sheet, placed after it.
$source.Copy($null, $source)
$sheet3 = $wb.ActiveSheet
$sheet3.Name = "created_sheet_20230901_143340"

# --- Second new sheet: created_sheet_20230901_144622 -----------------------
# Another duplicate, placed after the one we just created.
$source.Copy($null, $sheet3)
$sheet4 = $wb.ActiveSheet
$sheet4.Name = "created_sheet_20230901_144622"

# Copy column A values (A2:A5) and paste them into column C (C2:C5).
$sheet4.Range("A2:A5").Copy($sheet4.Range("C2:C5"))

Write-Host "Sheets now:"
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
